$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2302266.5
$ws.Range("J17").Value = 2302266.5
$ws.Range("L17").Value = 6906799.5
$ws.Range("N17").Value = -6907135.5
$ws.Range("H70").Value = 1439.0555
$ws.Range("I70").Value = 1066.6666
$ws.Range("J70").Value = 1513.5333
$ws.Range("K70").Value = 3199.9998
$ws.Range("L70").Value = 4540.5999
$ws.Range("M70").Value = -2929.9998
$ws.Range("N70").Value = -5080.5999
$ws.Range("H73").Value = 1439.0555
$ws.Range("I73").Value = 1066.6666
$ws.Range("J73").Value = 1513.5333
$ws.Range("K73").Value = 3199.9998
$ws.Range("L73").Value = 4540.5999
$ws.Range("M73").Value = -2263.9998
$ws.Range("N73").Value = -6412.5999
$ws.Range("H100").Value = 2003
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459
$ws.Range("H137").Value = 7408125
$ws.Range("I137").Value = 655.43475
$ws.Range("K137").Value = 1966.30425
$ws.Range("M137").Value = 583.6957499999999
$ws.Range("H138").Value = 6174303
$ws.Range("J138").Value = 2512.182
$ws.Range("L138").Value = 7536.545999999999
$ws.Range("N138").Value = -17816.546

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1549.9333
$ws.Range("I45").Value = 1478.1052
$ws.Range("J45").Value = 1674
$ws.Range("K45").Value = 1478.1052
$ws.Range("L45").Value = 1674
$ws.Range("M45").Value = -1101.1052
$ws.Range("N45").Value = -2428
$ws.Range("H61").Value = 8197900.5
$ws.Range("I61").Value = 10001193
$ws.Range("J61").Value = 1113.091
$ws.Range("K61").Value = 10001193
$ws.Range("L61").Value = 1113.091
$ws.Range("M61").Value = -10000981
$ws.Range("N61").Value = -1537.091
$ws.Range("H74").Value = 13515323
$ws.Range("I74").Value = 16130633
$ws.Range("J74").Value = 2885.6667
$ws.Range("K74").Value = 16130633
$ws.Range("L74").Value = 2885.6667
$ws.Range("M74").Value = -16129759
$ws.Range("N74").Value = -4633.6667
$ws.Range("H77").Value = 13515323
$ws.Range("I77").Value = 16130633
$ws.Range("J77").Value = 2885.6667
$ws.Range("K77").Value = 80653165
$ws.Range("L77").Value = 14428.3335
$ws.Range("M77").Value = -80648797
$ws.Range("N77").Value = -23164.3335
$ws.Range("H136").Value = 8197900.5
$ws.Range("I136").Value = 10001193
$ws.Range("J136").Value = 1113.091
$ws.Range("K136").Value = 30003579
$ws.Range("L136").Value = 3339.273
$ws.Range("M136").Value = -30001029
$ws.Range("N136").Value = -8439.272999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 689.439
$ws.Range("I94").Value = 576.30554
$ws.Range("J94").Value = 1504
$ws.Range("K94").Value = 576.30554
$ws.Range("L94").Value = 1504
$ws.Range("M94").Value = -125.30554
$ws.Range("N94").Value = -2406
$ws.Range("H134").Value = 3181.2263
$ws.Range("I134").Value = 2429.8684
$ws.Range("J134").Value = 5084.6665
$ws.Range("K134").Value = 7289.6052
$ws.Range("L134").Value = 15253.9995
$ws.Range("M134").Value = -4754.6052
$ws.Range("N134").Value = -20323.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4903800.5
$ws.Range("I31").Value = 1454.8182
$ws.Range("J31").Value = 25644494
$ws.Range("K31").Value = 1454.8182
$ws.Range("L31").Value = 25644494
$ws.Range("M31").Value = -1159.8182
$ws.Range("N31").Value = -25645084
$ws.Range("H34").Value = 4903800.5
$ws.Range("I34").Value = 1454.8182
$ws.Range("J34").Value = 25644494
$ws.Range("K34").Value = 1454.8182
$ws.Range("L34").Value = 25644494
$ws.Range("M34").Value = -1252.8182
$ws.Range("N34").Value = -25644898
$ws.Range("H58").Value = 1125.1464
$ws.Range("I58").Value = 499.42426
$ws.Range("J58").Value = 3706.25
$ws.Range("K58").Value = 499.42426
$ws.Range("L58").Value = 3706.25
$ws.Range("M58").Value = -296.42426
$ws.Range("N58").Value = -4112.25
$ws.Range("H124").Value = 17326
$ws.Range("J124").Value = 17326
$ws.Range("L124").Value = 17326
$ws.Range("N124").Value = -22236
$ws.Range("H132").Value = 8476223
$ws.Range("I132").Value = 10870878
$ws.Range("J132").Value = 2831.6924
$ws.Range("K132").Value = 32612634
$ws.Range("L132").Value = 8495.0772
$ws.Range("M132").Value = -32610104
$ws.Range("N132").Value = -13555.0772
$ws.Range("H134").Value = 1387.6222
$ws.Range("I134").Value = 1238.2094
$ws.Range("J134").Value = 4600
$ws.Range("K134").Value = 3714.6282
$ws.Range("L134").Value = 13800
$ws.Range("M134").Value = -1179.6282
$ws.Range("N134").Value = -18870
$ws.Range("H136").Value = 1125.1464
$ws.Range("I136").Value = 499.42426
$ws.Range("J136").Value = 3706.25
$ws.Range("K136").Value = 1498.27278
$ws.Range("L136").Value = 11118.75
$ws.Range("M136").Value = 1051.72722
$ws.Range("N136").Value = -16218.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 840.2353000000001
$ws.Range("I121").Value = 258.33334
$ws.Range("J121").Value = 1157.6364
$ws.Range("K121").Value = 775.0000200000001
$ws.Range("L121").Value = 3472.9092
$ws.Range("M121").Value = 534.9999799999999
$ws.Range("N121").Value = -6092.9092

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11497481
$ws.Range("I80").Value = 20835690
$ws.Range("K80").Value = 20835690
$ws.Range("M80").Value = -20834692
$ws.Range("H83").Value = 11497481
$ws.Range("I83").Value = 20835690
$ws.Range("K83").Value = 104178450
$ws.Range("M83").Value = -104173458
$ws.Range("H126").Value = 3468.3872
$ws.Range("I126").Value = 2165.4
$ws.Range("J126").Value = 4088.8572
$ws.Range("K126").Value = 6496.200000000001
$ws.Range("L126").Value = 12266.5716
$ws.Range("M126").Value = -4026.200000000001
$ws.Range("N126").Value = -17206.5716

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5560422.5
$ws.Range("I132").Value = 2704.2837
$ws.Range("K132").Value = 8112.8511
$ws.Range("M132").Value = -5582.8511
$ws.Range("H136").Value = 12824387
$ws.Range("I136").Value = 17242550
$ws.Range("J136").Value = 11711
$ws.Range("K136").Value = 51727650
$ws.Range("L136").Value = 35133
$ws.Range("M136").Value = -51725100
$ws.Range("N136").Value = -40233

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3734.2104
$ws.Range("I96").Value = 3186.3635
$ws.Range("J96").Value = 4487.5
$ws.Range("K96").Value = 3186.3635
$ws.Range("L96").Value = 4487.5
$ws.Range("M96").Value = -1813.3635
$ws.Range("N96").Value = -7233.5
$ws.Range("H132").Value = 1497.2858
$ws.Range("I132").Value = 1129.5238
$ws.Range("J132").Value = 2600.5715
$ws.Range("K132").Value = 3388.5714
$ws.Range("L132").Value = 7801.7145
$ws.Range("M132").Value = -858.5713999999998
$ws.Range("N132").Value = -12861.7145
$ws.Range("H136").Value = 1075.5
$ws.Range("I136").Value = 842.2
$ws.Range("J136").Value = 2242
$ws.Range("K136").Value = 2526.6
$ws.Range("L136").Value = 6726
$ws.Range("M136").Value = 23.39999999999964
$ws.Range("N136").Value = -11826
